$d = $word.ActiveDocument

# Locate the "1000" in "1000±100 ft to win" (the target altitude value).
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("1000", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find '1000' in the document"
}

$start = $find.Parent.Start

# Insert an extra "0" between the first two and last two digits, turning
# "1000" into "10000" (i.e. 10000 ft +/- 100 ft).
$insertPoint = $d.Range($start + 2, $start + 2)
$insertPoint.InsertBefore("0")

# Nudge formatting on the newly-typed character (as real Word's live
# proofing/formatting pass would) so it stays its own run instead of
# silently re-merging with its neighbours.
$toggleRange = $d.Range($start + 2, $start + 3)
$toggleRange.Font.Bold = 1
$toggleRange.Font.Bold = 0
